$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2023-01-06 12:56:11"
$newTimestamp = "2023-01-06 20:49:25"

# Update the timestamp column (O) for every data row (2 through 398).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Row 14: Buttertoast now shows "Online kein Bestand" (out of stock online).
$ws.Range("M14").Value = "Buttertoast 10 Scheiben - Online kein Bestand 2.20 Schweizer Franken"

# Row 21: ratingAmount bumped from 7 to 8.
$ws.Range("D21").Value = 8

# Row 40: Oelz Premium Drei-Korn Toast now out of stock online.
$ws.Range("M40").Value = "Ölz Premium Drei-Korn Toast - Online kein Bestand 3.40 Schweizer Franken"

# Row 52: Naturaplan Bio Mehrkorntoast mit Dinkel now out of stock online.
$ws.Range("M52").Value = "Naturaplan Bio Mehrkorntoast mit Dinkel - Online kein Bestand 2.95 Schweizer Franken"

# Row 75: ratingAmount bumped from 16 to 17, ratingValue dropped from 4.5 to 4.
$ws.Range("D75").Value = 17
$ws.Range("E75").Value = 4

# Row 194: Mulino Bianco Weizen-Brot back in stock online (text reverted).
$ws.Range("M194").Value = "Mulino Bianco Weizen-Brot 2.30 Schweizer Franken"

# Row 203: ratingAmount bumped from 9 to 10, ratingValue raised from 3.5 to 4.
$ws.Range("D203").Value = 10
$ws.Range("E203").Value = 4
